# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K") holds per-start strikeout counts that
# were previously populated from a different stat ("Strike#"). This script
# rewrites column G with the corrected K values for each existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values, keyed by row number, replacing the old
# Strike#-derived values that were written into column G.
$kValues = @{
    2  = 2
    3  = 1
    4  = 4
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 4
    13 = 1
    14 = 3
    15 = 0
    16 = 3
    17 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
